$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

$ws.Range("D2").Value = 20417
$ws.Range("E2").Value = 735
$ws.Range("F2").Value = 735
$ws.Range("G2").Value = 973
$ws.Range("H2").Value = 645
$ws.Range("I2").Value = 650
$ws.Range("J2").Value = -5
$ws.Range("K2").Value = 22241
$ws.Range("L2").Value = 6707
$ws.Range("M2").Value = 15535
$ws.Range("N2").Value = 15347
$ws.Range("O2").Value = 188
$ws.Range("P2").Value = 304
$ws.Range("Q2").Value = 987
$ws.Range("R2").Value = -682
$ws.Range("S2").Value = -266
$ws.Range("T2").Value = 1059
$ws.Range("U2").Value = -72
$ws.Range("V2").Value = 954
$ws.Range("W2").Value = 3.6
$ws.Range("X2").Value = 3.16
$ws.Range("Y2").Value = 4.29
$ws.Range("Z2").Value = 2.9
$ws.Range("AA2").Value = 43.17
$ws.Range("AB2").Value = 5204.45
$ws.Range("AC2").Value = 10687
$ws.Range("AD2").Value = 23.63
$ws.Range("AE2").Value = 265392
$ws.Range("AF2").Value = 0.95
$ws.Range("AG2").Value = 4000
$ws.Range("AH2").Value = 1.58
$ws.Range("AI2").Value = 35.58
$ws.Range("AJ2").Value = 6082642
$ws.Range("D3").Value = 21816
$ws.Range("E3").Value = 1183
$ws.Range("F3").Value = 1183
$ws.Range("G3").Value = 1535
$ws.Range("H3").Value = 1174
$ws.Range("I3").Value = 1173
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 24187
$ws.Range("L3").Value = 7722
$ws.Range("M3").Value = 16465
$ws.Range("N3").Value = 16350
$ws.Range("O3").Value = 115
$ws.Range("P3").Value = 304
$ws.Range("Q3").Value = 2019
$ws.Range("R3").Value = -1677
$ws.Range("S3").Value = 183
$ws.Range("T3").Value = 1215
$ws.Range("U3").Value = 803
$ws.Range("V3").Value = 1521
$ws.Range("W3").Value = 5.42
$ws.Range("X3").Value = 5.38
$ws.Range("Y3").Value = 7.4
$ws.Range("Z3").Value = 5.06
$ws.Range("AA3").Value = 46.9
$ws.Range("AB3").Value = 5506.76
$ws.Range("AC3").Value = 19291
$ws.Range("AD3").Value = 22.78
$ws.Range("AE3").Value = 282744
$ws.Range("AF3").Value = 1.55
$ws.Range("AG3").Value = 4000
$ws.Range("AH3").Value = 0.91
$ws.Range("AI3").Value = 19.71
$ws.Range("AJ3").Value = 6082642
$ws.Range("D4").Value = 22170
$ws.Range("E4").Value = 897
$ws.Range("F4").Value = 897
$ws.Range("G4").Value = 2333
$ws.Range("H4").Value = 1992
$ws.Range("I4").Value = 1993
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 24813
$ws.Range("L4").Value = 6789
$ws.Range("M4").Value = 18024
$ws.Range("N4").Value = 17899
$ws.Range("O4").Value = 125
$ws.Range("P4").Value = 304
$ws.Range("Q4").Value = 2210
$ws.Range("R4").Value = -1357
$ws.Range("S4").Value = -760
$ws.Range("T4").Value = 574
$ws.Range("U4").Value = 1636
$ws.Range("V4").Value = 991
$ws.Range("W4").Value = 4.05
$ws.Range("X4").Value = 8.99
$ws.Range("Y4").Value = 11.64
$ws.Range("Z4").Value = 8.130000000000001
$ws.Range("AA4").Value = 37.67
$ws.Range("AB4").Value = 6058.8
$ws.Range("AC4").Value = 32764
$ws.Range("AD4").Value = 10.15
$ws.Range("AE4").Value = 309527
$ws.Range("AF4").Value = 1.07
$ws.Range("AG4").Value = 4000
$ws.Range("AH4").Value = 1.2
$ws.Range("AI4").Value = 11.61
$ws.Range("AJ4").Value = 6082642
$ws.Range("D5").Value = 22083
$ws.Range("E5").Value = 964
$ws.Range("F5").Value = 964
$ws.Range("G5").Value = 1194
$ws.Range("H5").Value = 907
$ws.Range("I5").Value = 907
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 24499
$ws.Range("L5").Value = 6088
$ws.Range("M5").Value = 18411
$ws.Range("N5").Value = 18284
$ws.Range("O5").Value = 127
$ws.Range("P5").Value = 304
$ws.Range("Q5").Value = 1034
$ws.Range("R5").Value = -792
$ws.Range("S5").Value = -470
$ws.Range("T5").Value = 1282
$ws.Range("U5").Value = -248
$ws.Range("V5").Value = 682
$ws.Range("W5").Value = 4.37
$ws.Range("X5").Value = 4.11
$ws.Range("Y5").Value = 5.01
$ws.Range("Z5").Value = 3.68
$ws.Range("AA5").Value = 33.07
$ws.Range("AB5").Value = 6274.87
$ws.Range("AC5").Value = 14905
$ws.Range("AD5").Value = 23.75
$ws.Range("AE5").Value = 316196
$ws.Range("AF5").Value = 1.12
$ws.Range("AG5").Value = 4000
$ws.Range("AH5").Value = 1.13
$ws.Range("AI5").Value = 25.51
$ws.Range("AJ5").Value = 6082642
$ws.Range("D6").Value = 22364
$ws.Range("E6").Value = 886
$ws.Range("F6").Value = 886
$ws.Range("G6").Value = 1123
$ws.Range("H6").Value = 843
$ws.Range("I6").Value = 843
$ws.Range("K6").Value = 25323
$ws.Range("L6").Value = 6345
$ws.Range("M6").Value = 18978
$ws.Range("N6").Value = 18851
$ws.Range("P6").Value = 304
$ws.Range("Q6").Value = 1469
$ws.Range("R6").Value = -1289
$ws.Range("S6").Value = -58
$ws.Range("T6").Value = 882
$ws.Range("U6").Value = 587
$ws.Range("V6").Value = 883
$ws.Range("W6").Value = 3.96
$ws.Range("X6").Value = 3.77
$ws.Range("Y6").Value = 4.54
$ws.Range("Z6").Value = 3.38
$ws.Range("AA6").Value = 33.44
$ws.Range("AB6").Value = 6454.97
$ws.Range("AC6").Value = 13858
$ws.Range("AD6").Value = 18.37
$ws.Range("AE6").Value = 325998
$ws.Range("AF6").Value = 0.78
$ws.Range("AG6").Value = 4000
$ws.Range("AH6").Value = 1.57
$ws.Range("AI6").Value = 27.44
$ws.Range("AJ6").Value = 6082642
$ws.Range("D7").Value = 23477
$ws.Range("E7").Value = 827
$ws.Range("G7").Value = 1072
$ws.Range("H7").Value = 752
$ws.Range("I7").Value = 751
$ws.Range("K7").Value = 26183
$ws.Range("L7").Value = 6712
$ws.Range("M7").Value = 19472
$ws.Range("N7").Value = 19351
$ws.Range("P7").Value = 302
$ws.Range("Q7").Value = 1505
$ws.Range("R7").Value = -1173
$ws.Range("S7").Value = -111
$ws.Range("T7").Value = 1048
$ws.Range("U7").Value = 439
$ws.Range("W7").Value = 3.52
$ws.Range("X7").Value = 3.2
$ws.Range("Y7").Value = 3.93
$ws.Range("Z7").Value = 2.92
$ws.Range("AA7").Value = 34.47
$ws.Range("AC7").Value = 12354
$ws.Range("AD7").Value = 18.7
$ws.Range("AE7").Value = 334633
$ws.Range("AF7").Value = 0.6899999999999999
$ws.Range("AG7").Value = 4000
$ws.Range("AH7").Value = 1.73
$ws.Range("AI7").Value = 32.38
$ws.Range("D8").Value = 24747
$ws.Range("E8").Value = 963
$ws.Range("G8").Value = 1225
$ws.Range("H8").Value = 924
$ws.Range("I8").Value = 923
$ws.Range("K8").Value = 27071
$ws.Range("L8").Value = 6953
$ws.Range("M8").Value = 20117
$ws.Range("N8").Value = 19990
$ws.Range("P8").Value = 302
$ws.Range("Q8").Value = 1638
$ws.Range("R8").Value = -1196
$ws.Range("S8").Value = -194
$ws.Range("T8").Value = 950
$ws.Range("U8").Value = 570
$ws.Range("W8").Value = 3.89
$ws.Range("X8").Value = 3.73
$ws.Range("Y8").Value = 4.69
$ws.Range("Z8").Value = 3.47
$ws.Range("AA8").Value = 34.56
$ws.Range("AC8").Value = 15174
$ws.Range("AD8").Value = 15.22
$ws.Range("AE8").Value = 345686
$ws.Range("AF8").Value = 0.67
$ws.Range("AG8").Value = 4000
$ws.Range("AH8").Value = 1.73
$ws.Range("AI8").Value = 26.36
$ws.Range("D9").Value = 25946
$ws.Range("E9").Value = 1072
$ws.Range("G9").Value = 1337
$ws.Range("H9").Value = 1008
$ws.Range("I9").Value = 1008
$ws.Range("K9").Value = 28020
$ws.Range("L9").Value = 7178
$ws.Range("M9").Value = 20843
$ws.Range("N9").Value = 20710
$ws.Range("P9").Value = 302
$ws.Range("Q9").Value = 1728
$ws.Range("R9").Value = -1253
$ws.Range("S9").Value = -162
$ws.Range("T9").Value = 970
$ws.Range("U9").Value = 634
$ws.Range("W9").Value = 4.13
$ws.Range("X9").Value = 3.89
$ws.Range("Y9").Value = 4.95
$ws.Range("Z9").Value = 3.66
$ws.Range("AA9").Value = 34.44
$ws.Range("AC9").Value = 16570
$ws.Range("AD9").Value = 13.94
$ws.Range("AE9").Value = 358149
$ws.Range("AF9").Value = 0.64
$ws.Range("AG9").Value = 4091
$ws.Range("AH9").Value = 1.77
$ws.Range("AI9").Value = 24.69
